$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.816235780715942
$ws.Range("B1").Value = 6.017861843109131
$ws.Range("C1").Value = 5.078186511993408
$ws.Range("D1").Value = 5.908387184143066
$ws.Range("E1").Value = 4.096496105194092
